$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.320.17"
$ws.Range("E2").Value = "  +0.58%  "
$ws.Range("D3").Value = "1.872.49"
$ws.Range("E3").Value = "  +0.62%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Formula = "'0.7120"
$ws.Range("E5").Value = "  +0.56%  "
$ws.Range("D6").Formula = "'241.56"
$ws.Range("E6").Value = "  +0.32%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Formula = "'0.07788"
$ws.Range("E8").Value = "  +1.92%  "
$ws.Range("D9").Formula = "'0.3110"
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").Formula = "'25.08"
$ws.Range("E10").Value = "  +1.82%  "
$ws.Range("D11").Formula = "'0.08403"
$ws.Range("E11").Value = "  +0.94%  "
$ws.Range("D12").Value = "1.864.57"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("D13").Formula = "'5.233"
$ws.Range("E13").Value = "  +1.13%  "
$ws.Range("D14").Formula = "'0.7115"
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("D15").Formula = "'91.06"
$ws.Range("E15").Value = "  -0.05%  "
$ws.Range("D16").Value = "29.329.55"
$ws.Range("E16").Value = "  +0.57%  "
$ws.Range("D17").Formula = "'6.085"
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("D18").Formula = "'0.000008221"
$ws.Range("E18").Value = "  +5.28%  "
$ws.Range("D19").Formula = "'240.27"
$ws.Range("E19").Value = "  -0.85%  "
$ws.Range("D20").Formula = "'13.20"
$ws.Range("E20").Value = "  +1.06%  "
$ws.Range("D21").Value = "2.121.89"
$ws.Range("E21").Value = "  +0.29%  "
$ws.Range("D22").Formula = "'0.9999"
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Formula = "'7.762"
$ws.Range("E23").Value = "  -0.91%  "
$ws.Range("D24").Formula = "'1.000"
$ws.Range("E25").Value = "  +0.69%  "
$ws.Range("D26").Formula = "'162.87"
$ws.Range("E26").Value = "  -0.19%  "
$ws.Range("D27").Formula = "'9.026"
$ws.Range("E27").Value = "  +1.11%  "
$ws.Range("D28").Formula = "'18.50"
$ws.Range("E28").Value = "  +0.30%  "
$ws.Range("E29").Value = "  +0.86%  "
$ws.Range("D30").Formula = "'4.416"
$ws.Range("E30").Value = "  +0.48%  "
$ws.Range("E31").Value = "  -1.90%  "
$ws.Range("D32").Formula = "'4.308"
$ws.Range("E32").Value = "  +2.54%  "
$ws.Range("E33").Value = "  +3.30%  "
$ws.Range("D34").Formula = "'1.935"
$ws.Range("D35").Formula = "'1.176"
$ws.Range("E35").Value = "  +1.31%  "
$ws.Range("D36").Formula = "'0.7448"
$ws.Range("E36").Value = "  -6.58%  "
$ws.Range("D37").Formula = "'2.699"
$ws.Range("E37").Value = "  +0.71%  "
$ws.Range("D38").Formula = "'0.01870"
$ws.Range("E38").Value = "  +1.67%  "
$ws.Range("D39").Value = "1.225.91"
$ws.Range("E39").Value = "  +5.23%  "
$ws.Range("D40").Formula = "'2.727"
$ws.Range("E40").Value = "  +1.09%  "
$ws.Range("D41").Formula = "'6.560"
$ws.Range("E41").Value = "  +6.33%  "
$ws.Range("D42").Formula = "'110.49"
$ws.Range("E42").Value = "  +8.22%  "
$ws.Range("D43").Formula = "'0.8876"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("D44").Formula = "'72.62"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").Formula = "'1.000"
$ws.Range("E45").Value = "  +0.02%  "
$ws.Range("D46").Value = "2.019.38"
$ws.Range("E46").Value = "  +0.49%  "
$ws.Range("E47").Value = "  +1.94%  "
$ws.Range("D48").Formula = "'0.5198"
$ws.Range("E48").Value = "  +0.35%  "
$ws.Range("E49").Value = "  +2.00%  "
$ws.Range("D50").Formula = "'9.405"
$ws.Range("E50").Value = "  +0.98%  "
$ws.Range("D51").Formula = "'0.4316"
$ws.Range("E51").Value = "  +1.19%  "
